$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(2, 'Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '25.599.67', '  -5.97%  '),
    @(3, 'Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '1.806.44', '  -5.20%  '),
    @(4, 'TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '0.9999', '  -0.06%  '),
    @(5, 'BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '276.66', '  -9.65%  '),
    @(6, 'USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '0.9994', '  -0.11%  '),
    @(7, 'XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.5036', '  -6.27%  '),
    @(8, 'Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.3504', '  -7.98%  '),
    @(9, 'OKB', 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb', '43.63', '  -5.15%  '),
    @(10, 'Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.06668', '  -8.55%  '),
    @(11, 'Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '20.01', '  -9.96%  '),
    @(12, 'Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '0.8357', '  -7.68%  '),
    @(13, 'TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.07806', '  -4.78%  '),
    @(14, 'WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '1.812.04', '  +42.62%  '),
    @(15, 'Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '5.058', '  -5.30%  '),
    @(16, 'Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '87.35', '  -8.79%  '),
    @(17, 'BinanceUSD', 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd', '0.9984', '  -0.08%  '),
    @(18, 'Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '13.90', '  -6.43%  '),
    @(19, 'Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '1.000', '  +0.01%  '),
    @(20, 'ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.000007934', '  -8.46%  '),
    @(21, 'WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '25.682.19', '  -5.80%  '),
    @(22, 'Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '4.709', '  -6.70%  '),
    @(23, 'WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '2.038.10', '  +67.24%  '),
    @(24, 'Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '9.972', '  -7.58%  '),
    @(25, 'Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '6.050', '  -7.24%  '),
    @(26, 'Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '141.41', '  -4.89%  '),
    @(27, 'LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '2.117', '  -8.34%  '),
    @(28, 'Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '1.652', '  -5.47%  '),
    @(29, 'EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '16.99', '  -7.63%  '),
    @(30, 'BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '108.74', '  -6.81%  '),
    @(31, 'InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '4.298', '  -11.28%  '),
    @(32, 'Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '4.216', '  -10.82%  '),
    @(33, 'Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.08890', '  -3.56%  '),
    @(34, 'Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.04801', '  -5.48%  '),
    @(35, 'ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '0.7315', '  -11.65%  '),
    @(36, 'ARBITRUM', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', '1.126', '  -7.62%  '),
    @(37, 'HuobiToken', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht', '2.840', '  -5.36%  '),
    @(38, 'Frax', 'https://coinranking.com/coin/KfWtaeV1W+frax-frax', '0.9988', '  -0.08%  '),
    @(39, 'MXToken', 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx', '3.033', '  -8.57%  '),
    @(40, 'VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.01862', '  -7.02%  '),
    @(41, 'TheSandbox', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand', '0.5175', '  -13.11%  '),
    @(42, 'RenderToken', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', '2.297', '  -14.01%  '),
    @(43, 'TrustWalletToken', 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt', '0.9617', '  -10.88%  '),
    @(44, 'Quant', 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt', '113.64', '  -2.50%  '),
    @(45, 'FraxShare', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', '6.185', '  -7.24%  '),
    @(46, 'Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '8.064', '  -13.53%  '),
    @(47, 'PaxDollar', 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp', '0.9985', '  -0.17%  '),
    @(48, 'Decentraland', 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana', '0.4577', '  -11.34%  '),
    @(49, 'Algorand', 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo', '0.1383', '  -9.57%  '),
    @(50, 'EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '9.209', '  -9.93%  '),
    @(51, 'Elrond', 'https://coinranking.com/coin/omwkOTglq+elrond-egld', '35.73', '  -6.90%  ')
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = "'" + $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}
